$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
$ws.Cells.Item(1,1).Value = "nama_survei"
$ws.Cells.Item(1,2).Value = "lokasi_survei"
$ws.Cells.Item(1,3).Value = "kode_desa"
$ws.Cells.Item(1,4).Value = "kode_kecamatan"
$ws.Cells.Item(1,5).Value = "kode_kabupaten"
$ws.Cells.Item(1,6).Value = "kode_provinsi"
$ws.Cells.Item(1,7).Value = "kro"
$ws.Cells.Item(1,8).Value = "jadwal"
$ws.Cells.Item(1,9).Value = "tim"

# --- Row 2 (sample/example data) ---
$ws.Cells.Item(2,1).Value = "contoh"
$ws.Cells.Item(2,2).Value = "contoh"
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(2,4).Value = 1
$ws.Cells.Item(2,5).ClearContents()
$ws.Cells.Item(2,6).ClearContents()
$ws.Cells.Item(2,7).Value = "contoh"
$ws.Cells.Item(2,8).Value = "01-12-2029"
$ws.Cells.Item(2,9).Value = "contoh"

# --- Selection moves from N3 to M4 ---
$ws.Range("M4").Select()
